$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are formatted as Text so numeric-looking / percent-looking
# strings are stored verbatim (matching the source data, which are literal
# text labels, not real numbers/percentages).
$cells = @(
    "D2", "E2", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10",
    "D11", "E11", "D12", "E12", "D13", "E13", "D14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18",
    "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26",
    "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46",
    "E46", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (as text).
$ws.Range("D2").Value = '303.94'
$ws.Range("E2").Value = '-0.03%'
$ws.Range("E3").Value = '3.18%'
$ws.Range("D4").Value = '5.010'
$ws.Range("E4").Value = '-1.70%'
$ws.Range("D5").Value = '0.07769'
$ws.Range("E5").Value = '-0.37%'
$ws.Range("D6").Value = '2.103'
$ws.Range("E6").Value = '-7.56%'
$ws.Range("D7").Value = '8.031'
$ws.Range("E7").Value = '-0.99%'
$ws.Range("D8").Value = '0.9209'
$ws.Range("E8").Value = '-0.71%'
$ws.Range("D9").Value = '0.09869'
$ws.Range("E9").Value = '1.41%'
$ws.Range("D10").Value = '0.1864'
$ws.Range("E10").Value = '2.14%'
$ws.Range("D11").Value = '0.08654'
$ws.Range("E11").Value = '-0.28%'
$ws.Range("D12").Value = '0.03591'
$ws.Range("E12").Value = '5.03%'
$ws.Range("D13").Value = '0.09976'
$ws.Range("E13").Value = '0.36%'
$ws.Range("D14").Value = '0.001480'
$ws.Range("D15").Value = '0.005691'
$ws.Range("E15").Value = '-0.76%'
$ws.Range("D16").Value = '3.460'
$ws.Range("E16").Value = '-0.50%'
$ws.Range("D17").Value = '4.052'
$ws.Range("E17").Value = '0.50%'
$ws.Range("D18").Value = '2.533'
$ws.Range("E18").Value = '18.00%'
$ws.Range("E19").Value = '-0.39%'
$ws.Range("D20").Value = '0.1332'
$ws.Range("E20").Value = '0.71%'
$ws.Range("D21").Value = '4.945'
$ws.Range("E21").Value = '8.38%'
$ws.Range("D22").Value = '0.2210'
$ws.Range("E22").Value = '-1.21%'
$ws.Range("D23").Value = '0.04605'
$ws.Range("E23").Value = '-1.76%'
$ws.Range("D24").Value = '0.005137'
$ws.Range("E24").Value = '13.15%'
$ws.Range("D25").Value = '0.001237'
$ws.Range("E25").Value = '-0.45%'
$ws.Range("D26").Value = '0.0001409'
$ws.Range("E26").Value = '8.19%'
$ws.Range("D39").Value = '0.01791'
$ws.Range("E39").Value = '2.21%'
$ws.Range("D40").Value = '0.04666'
$ws.Range("E40").Value = '-0.89%'
$ws.Range("D41").Value = '0.007711'
$ws.Range("E41").Value = '-3.38%'
$ws.Range("D42").Value = '0.1397'
$ws.Range("E42").Value = '-1.64%'
$ws.Range("D43").Value = '0.007629'
$ws.Range("E43").Value = '-4.77%'
$ws.Range("D44").Value = '0.002224'
$ws.Range("E44").Value = '-3.04%'
$ws.Range("D45").Value = '0.01042'
$ws.Range("E45").Value = '14.44%'
$ws.Range("D46").Value = '0.00006320'
$ws.Range("E46").Value = '1.29%'
$ws.Range("E47").Value = '0.33%'
$ws.Range("D48").Value = '0.0005827'
$ws.Range("E48").Value = '0.46%'
$ws.Range("D49").Value = '33.86'
$ws.Range("E49").Value = '497.74%'
$ws.Range("D50").Value = '0.002009'
$ws.Range("E50").Value = '-25.43%'
$ws.Range("D51").Value = '0.00002110'
$ws.Range("E51").Value = '0.33%'
